# Apply attendance marks to Sheet1 of the workbook.
# Columns: A=Date, B=Roll, C=Name, D=Total Attendance Count, E=Real,
#          F=Duplicate, G=Invalid, H=Absent
# For each date row (3..18), set the appropriate cells from 0 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where attendance was "Real": set D (Total Attendance Count) and E (Real) to 1
$realRows = @(5, 6, 10, 11, 12, 14, 15)
foreach ($r in $realRows) {
    $ws.Cells.Item($r, 4).Value = 1   # D - Total Attendance Count
    $ws.Cells.Item($r, 5).Value = 1   # E - Real
}

# Rows where attendance was "Invalid": set G (Invalid) and H (Absent) to 1
$invalidRows = @(3, 17)
foreach ($r in $invalidRows) {
    $ws.Cells.Item($r, 7).Value = 1   # G - Invalid
    $ws.Cells.Item($r, 8).Value = 1   # H - Absent
}

# Rows where the student was simply "Absent": set H (Absent) to 1
$absentRows = @(4, 7, 8, 9, 13, 16, 18)
foreach ($r in $absentRows) {
    $ws.Cells.Item($r, 8).Value = 1   # H - Absent
}
